# Update cryptocurrency price and volume data in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so that numeric-looking
# strings (e.g. "1.000", "4.830") retain their exact original formatting
# instead of being auto-converted into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.684.57"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "1.690.75"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "218.09"
$ws.Range("D6").Value = "0.5344"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.2687"
$ws.Range("E8").Value = "  +5.02%  "
$ws.Range("D9").Value = "0.06446"
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("D10").Value = "21.78"
$ws.Range("E10").Value = "  +7.59%  "
$ws.Range("D11").Value = "0.07804"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").Value = "1.693.31"
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("D13").Value = "4.515"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").Value = "0.5663"
$ws.Range("E14").Value = "  +4.70%  "
$ws.Range("D15").Value = "0.0₅8483"
$ws.Range("E15").Value = "  +7.40%  "
$ws.Range("D16").Value = "66.63"
$ws.Range("E16").Value = "  +3.46%  "
$ws.Range("D17").Value = "26.679.38"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "4.830"
$ws.Range("E19").Value = "  +4.53%  "
$ws.Range("D20").Value = "196.23"
$ws.Range("E20").Value = "  +6.72%  "
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  +4.57%  "
$ws.Range("D22").Value = "6.411"
$ws.Range("E22").Value = "  +5.57%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").Value = "143.90"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "0.1287"
$ws.Range("E25").Value = "  +7.37%  "
$ws.Range("D26").Value = "7.499"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").Value = "16.32"
$ws.Range("E27").Value = "  +5.68%  "
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("D29").Value = "0.06203"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").Value = "3.618"
$ws.Range("E31").Value = "  +8.54%  "
$ws.Range("D32").Value = "3.482"
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("D33").Value = "1.706"
$ws.Range("E33").Value = "  +6.44%  "
$ws.Range("D34").Value = "1.018"
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("D35").Value = "2.798"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").Value = "2.421"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").Value = "0.5749"
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("D38").Value = "0.01655"
$ws.Range("E38").Value = "  +4.11%  "
$ws.Range("D39").Value = "6.009"
$ws.Range("E39").Value = "  +6.39%  "
$ws.Range("D40").Value = "1.082.18"
$ws.Range("E40").Value = "  +5.52%  "
$ws.Range("D41").Value = "0.8671"
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("D43").Value = "100.56"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").Value = "1.836.77"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").Value = "0.0₈111"
$ws.Range("E45").Value = "  +5.73%  "
$ws.Range("D46").Value = "57.56"
$ws.Range("E46").Value = "  +6.28%  "
$ws.Range("D47").Value = "8.156"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").Value = "0.05223"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").Value = "6.108"
$ws.Range("E50").Value = "  +5.61%  "
$ws.Range("D51").Value = "0.4236"
$ws.Range("E51").Value = "  +0.19%  "
